# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Replaces the five mixed "CC / periodo 2506" rows (one per employee) with
# five rows for a single worker (ALDAIR ANTONIO HERNANDEZ PATERNINA,
# CC 1143364877) covering periods 2504-2508, updates the "Valor Mora" /
# "Salario Basico" amounts for the new periods, refreshes the summary
# totals (Cant. Trabajadores / Cant. Periodos / Valor Mora), and removes
# the three now-obsolete worker rows at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block -------------------------------------------------------
$ws.Range("E11").Value = 341520   # VALOR MORA (total)
$ws.Range("C13").Value = 1        # Cant. Trabajadores
$ws.Range("F13").Value = 5        # Cant. Periodos

# --- Detail rows (16-20): single worker across periods 2504-2508 --------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143364877"
$ws.Range("D16").Value = "ALDAIR ANTONIO HERNANDEZ PATERNINA"
$ws.Range("E16").Value = "2504"
$ws.Range("F16").Value = 68304
$ws.Range("G16").Value = 1707600

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143364877"
$ws.Range("D17").Value = "ALDAIR ANTONIO HERNANDEZ PATERNINA"
$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 68304
$ws.Range("G17").Value = 1707600

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143364877"
$ws.Range("D18").Value = "ALDAIR ANTONIO HERNANDEZ PATERNINA"
$ws.Range("E18").Value = "2506"
$ws.Range("F18").Value = 68304
$ws.Range("G18").Value = 1707600

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143364877"
$ws.Range("D19").Value = "ALDAIR ANTONIO HERNANDEZ PATERNINA"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 68304
$ws.Range("G19").Value = 1707600

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143364877"
$ws.Range("D20").Value = "ALDAIR ANTONIO HERNANDEZ PATERNINA"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 68304
$ws.Range("G20").Value = 1707600

# Row 20 is now the last row of the table, so it needs the table's
# "closing" bottom-border formatting that used to live on row 23 (the
# previous last row) before that row is removed.
$ws.Range("B23:J23").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)   # xlPasteFormats

# --- Remove the three obsolete worker rows (old rows 21-23) --------------
# Deleting shifts the trailing signature block (old rows 28-29) up to 25-26,
# matching the new used range B2:J26.
$ws.Range("B21:J23").EntireRow.Delete()
